$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the existing row 134: A134 timestamp gets a tiny precision correction,
# and E134/F134 (previously blank) get filled in with real values.
$ws.Range("A134").Value = 45986.43557956019
$ws.Range("A134").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E134").Value = 46.1
$ws.Range("F134").Value = 45986
$ws.Range("F134").NumberFormat = "YYYY-MM-DD"

# New rows 135-142: same product/pack/price pattern as the rows above them,
# each with its own timestamp, and price/date columns populated.
$newRows = @(
    @(135, 45986.48651076389, 45986),
    @(136, 45986.49729130787, 45986),
    @(137, 45987.42162665509, 45987),
    @(138, 45987.42379971065, 45987),
    @(139, 45987.42442510417, 45987),
    @(140, 45987.42470152778, 45987),
    @(141, 45987.42627121528, 45987),
    @(142, 45987.4267159375,  45987)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("A$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("B$r").Value = "EVOWHEY PROTEIN"
    $ws.Range("C$r").Value = "Pack (5x500g)"
    $ws.Range("D$r").Value = "46,10€"
    $ws.Range("E$r").Value = 46.1
    $ws.Range("F$r").Value = $row[2]
    $ws.Range("F$r").NumberFormat = "YYYY-MM-DD"
}

# Final new row 143: same product/pack/price text, but precio_num and
# fecha_dia are left blank (mirrors how row 134 looked before this edit).
$ws.Range("A143").Value = 45987.42767409614
$ws.Range("A143").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B143").Value = "EVOWHEY PROTEIN"
$ws.Range("C143").Value = "Pack (5x500g)"
$ws.Range("D143").Value = "46,10€"
